$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 176.84616
$ws.Range("I5").Value = 176.84616
$ws.Range("K5").Value = 176.84616
$ws.Range("M5").Value = -61.84616
$ws.Range("H6").Value = 115
$ws.Range("I6").Value = 115
$ws.Range("K6").Value = 345
$ws.Range("M6").Value = -233
$ws.Range("H38").Value = 2867.1428
$ws.Range("I38").Value = 37.5
$ws.Range("J38").Value = 3999
$ws.Range("K38").Value = 112.5
$ws.Range("L38").Value = 11997
$ws.Range("M38").Value = 259.5
$ws.Range("N38").Value = -12741
$ws.Range("H53").Value = 6140.4287
$ws.Range("I53").Value = 997.2857
$ws.Range("J53").Value = 16426.715
$ws.Range("K53").Value = 997.2857
$ws.Range("L53").Value = 16426.715
$ws.Range("M53").Value = -360.2857
$ws.Range("N53").Value = -17700.715
$ws.Range("H88").Value = 1000
$ws.Range("I88").Value = 1000
$ws.Range("K88").Value = 1000
$ws.Range("M88").Value = -594
$ws.Range("H91").Value = 1000
$ws.Range("I91").Value = 1000
$ws.Range("K91").Value = 1000
$ws.Range("M91").Value = 404
$ws.Range("H113").Value = 12005.267
$ws.Range("I113").Value = 15340.286
$ws.Range("J113").Value = 9087.125
$ws.Range("K113").Value = 15340.286
$ws.Range("L113").Value = 9087.125
$ws.Range("M113").Value = -12086.286
$ws.Range("N113").Value = -15595.125

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 4324.353
$ws.Range("I2").Value = 924.3333
$ws.Range("J2").Value = 8149.375
$ws.Range("K2").Value = 924.3333
$ws.Range("L2").Value = 8149.375
$ws.Range("M2").Value = -811.3333
$ws.Range("N2").Value = -8375.375
$ws.Range("H32").Value = 9505.360000000001
$ws.Range("I32").Value = 3564.7273
$ws.Range("K32").Value = 3564.7273
$ws.Range("M32").Value = -3277.7273
$ws.Range("H45").Value = 1578.6666
$ws.Range("I45").Value = 1427.2
$ws.Range("K45").Value = 1427.2
$ws.Range("M45").Value = -1050.2
$ws.Range("H116").Value = 4324.353
$ws.Range("I116").Value = 924.3333
$ws.Range("J116").Value = 8149.375
$ws.Range("K116").Value = 924.3333
$ws.Range("L116").Value = 8149.375
$ws.Range("M116").Value = 1369.6667
$ws.Range("N116").Value = -12737.375

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 4324.353
$ws.Range("I3").Value = 924.3333
$ws.Range("J3").Value = 8149.375
$ws.Range("K3").Value = 924.3333
$ws.Range("L3").Value = 8149.375
$ws.Range("M3").Value = -810.3333
$ws.Range("N3").Value = -8377.375
$ws.Range("H86").Value = 1353.5
$ws.Range("J86").Value = 1353.5
$ws.Range("L86").Value = 1353.5
$ws.Range("N86").Value = -3599.5
$ws.Range("H89").Value = 1353.5
$ws.Range("J89").Value = 1353.5
$ws.Range("L89").Value = 6767.5
$ws.Range("N89").Value = -17999.5
$ws.Range("I134").Value = 2031.6
$ws.Range("J134").Value = 22621.455
$ws.Range("K134").Value = 6094.799999999999
$ws.Range("L134").Value = 67864.36500000001
$ws.Range("M134").Value = -3559.799999999999
$ws.Range("N134").Value = -72934.36500000001

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 35887.535
$ws.Range("I31").Value = 35499.668
$ws.Range("J31").Value = 35984.5
$ws.Range("K31").Value = 35499.668
$ws.Range("L31").Value = 35984.5
$ws.Range("M31").Value = -35204.668
$ws.Range("N31").Value = -36574.5
$ws.Range("H34").Value = 35887.535
$ws.Range("I34").Value = 35499.668
$ws.Range("J34").Value = 35984.5
$ws.Range("K34").Value = 35499.668
$ws.Range("L34").Value = 35984.5
$ws.Range("M34").Value = -35297.668
$ws.Range("N34").Value = -36388.5
$ws.Range("H99").Value = 11765.25
$ws.Range("I99").Value = 3534
$ws.Range("K99").Value = 3534
$ws.Range("M99").Value = -2036
$ws.Range("H126").Value = 11765.25
$ws.Range("I126").Value = 3534
$ws.Range("K126").Value = 10602
$ws.Range("M126").Value = -8132
$ws.Range("H132").Value = 10046.294
$ws.Range("I132").Value = 2727.3
$ws.Range("K132").Value = 8181.900000000001
$ws.Range("M132").Value = -5651.900000000001
$ws.Range("H134").Value = 41675492
$ws.Range("I134").Value = 2111.5557
$ws.Range("J134").Value = 66679520
$ws.Range("K134").Value = 6334.6671
$ws.Range("L134").Value = 200038560
$ws.Range("M134").Value = -3799.6671
$ws.Range("N134").Value = -200043630

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1743410.4
$ws.Range("J5").Value = 3485514.5
$ws.Range("L5").Value = 10456543.5
$ws.Range("N5").Value = -10456767.5
$ws.Range("H12").Value = 128.5
$ws.Range("I12").Value = 152
$ws.Range("J12").Value = 120.666664
$ws.Range("K12").Value = 456
$ws.Range("L12").Value = 361.999992
$ws.Range("M12").Value = -283
$ws.Range("N12").Value = -707.999992
$ws.Range("H34").Value = 1368.6666
$ws.Range("I34").Value = 1869.8334
$ws.Range("J34").Value = 366.33334
$ws.Range("K34").Value = 5609.5002
$ws.Range("L34").Value = 1099.00002
$ws.Range("M34").Value = -5525.5002
$ws.Range("N34").Value = -1267.00002
$ws.Range("H98").Value = 7572.5
$ws.Range("J98").Value = 7572.5
$ws.Range("L98").Value = 22717.5
$ws.Range("N98").Value = -25713.5
$ws.Range("H135").Value = 1743410.4
$ws.Range("J135").Value = 3485514.5
$ws.Range("L135").Value = 31369630.5
$ws.Range("N135").Value = -31374700.5

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H17").Value = 1575.5
$ws.Range("J17").Value = 1743.4286
$ws.Range("L17").Value = 1743.4286
$ws.Range("N17").Value = -2079.4286
$ws.Range("H132").Value = 14121.174
$ws.Range("I132").Value = 8613.65
$ws.Range("K132").Value = 25840.95
$ws.Range("M132").Value = -23310.95

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 3718.389
$ws.Range("I16").Value = 3718.389
$ws.Range("K16").Value = 3718.389
$ws.Range("M16").Value = -3548.389
$ws.Range("H22").Value = 4760.4
$ws.Range("I22").Value = 1780.875
$ws.Range("J22").Value = 6746.75
$ws.Range("K22").Value = 1780.875
$ws.Range("L22").Value = 6746.75
$ws.Range("M22").Value = -1485.875
$ws.Range("N22").Value = -7336.75
$ws.Range("H27").Value = 4760.4
$ws.Range("I27").Value = 1780.875
$ws.Range("J27").Value = 6746.75
$ws.Range("K27").Value = 1780.875
$ws.Range("L27").Value = 6746.75
$ws.Range("M27").Value = -1673.875
$ws.Range("N27").Value = -6960.75
$ws.Range("H55").Value = 2415.0967
$ws.Range("I55").Value = 1184.6666
$ws.Range("J55").Value = 3568.625
$ws.Range("K55").Value = 1184.6666
$ws.Range("L55").Value = 3568.625
$ws.Range("M55").Value = -1011.6666
$ws.Range("N55").Value = -3914.625
$ws.Range("H57").Value = 29500
$ws.Range("I57").Value = 29000
$ws.Range("J57").Value = 30000
$ws.Range("K57").Value = 29000
$ws.Range("L57").Value = 30000
$ws.Range("M57").Value = -28434
$ws.Range("N57").Value = -31132
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").Value = ""
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").Value = ""
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("N110").Value = ""
$ws.Range("H122").Value = 7252.2256
$ws.Range("I122").Value = 5173.6313
$ws.Range("J122").Value = 10543.333
$ws.Range("K122").Value = 15520.8939
$ws.Range("L122").Value = 31629.999
$ws.Range("M122").Value = -13070.8939
$ws.Range("N122").Value = -36529.999
$ws.Range("H133").Value = 70666
$ws.Range("I133").Value = 72000
$ws.Range("J133").Value = 69999
$ws.Range("K133").Value = 72000
$ws.Range("L133").Value = 69999
$ws.Range("M133").Value = -69470
$ws.Range("N133").Value = -75059
$ws.Range("H136").Value = 11736.152
$ws.Range("I136").Value = 10568.458
$ws.Range("J136").Value = 13010
$ws.Range("K136").Value = 31705.374
$ws.Range("L136").Value = 39030
$ws.Range("M136").Value = -29155.374
$ws.Range("N136").Value = -44130

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 6934.05
$ws.Range("I122").Value = 2835.182
$ws.Range("J122").Value = 11943.777
$ws.Range("K122").Value = 8505.545999999998
$ws.Range("L122").Value = 35831.331
$ws.Range("M122").Value = -6055.545999999998
$ws.Range("N122").Value = -40731.331
$ws.Range("H136").Value = 8289.543
$ws.Range("I136").Value = 1540.75
$ws.Range("K136").Value = 4622.25
$ws.Range("M136").Value = -2072.25
